# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" right after "总计" (becomes the 2nd tab),
#    pushing 2022-Q2 / 2022-Q1 / 2021-Q4 / 2021-Q3 / 2021-Q1 down by one slot.
# 2) Populate "2022-Q3" with the per-fund holdings for that quarter (same shape
#    as the other quarterly sheets: header row + 4 fund rows).
# 3) Update the "总计" (totals) sheet with a new top data row for 2022-Q3 and
#    re-number the existing rows beneath it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value as TEXT (keeps leading zeros / exact decimal strings
# like "004317" or "0.0519" instead of letting them be auto-coerced to
# numbers). We stash the format tweak and wipe it again afterwards via a
# format-only paste from an always-empty reference cell, so cells end up with
# the plain/default style - matching the source workbook's other sheets,
# where these text cells carry no explicit style index.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# =====================================================================
# Step 1: insert the new "2022-Q3" worksheet right after "总计"
# =====================================================================
$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"

# =====================================================================
# Step 2: fill in the 2022-Q3 holdings sheet
# =====================================================================
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Header row (B1:H1) + style-only A-column header cell (A1 has no value in
# the source sheets, just copy the real header labels straight across).
$q2Sheet.Range("B1:H1").Copy($q3Sheet.Range("B1:H1"))

# A-column "index" style (s=2) used on A2 in the other quarterly sheets.
$q2Sheet.Range("A2").Copy($q3Sheet.Range("A2"))
$q2Sheet.Range("A2").Copy($q3Sheet.Range("A3"))
$q2Sheet.Range("A2").Copy($q3Sheet.Range("A4"))
$q2Sheet.Range("A2").Copy($q3Sheet.Range("A5"))

$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("A4").Value = 2
$q3Sheet.Range("A5").Value = 3

Set-TextValue $q3Sheet.Range("B2") "004317"
Set-TextValue $q3Sheet.Range("C2") "前海开源沪港深裕鑫灵活配置混合C"
Set-TextValue $q3Sheet.Range("D2") "1.79"
Set-TextValue $q3Sheet.Range("E2") "70.17"
Set-TextValue $q3Sheet.Range("F2") "2.90"
Set-TextValue $q3Sheet.Range("G2") "0.0519"
$q3Sheet.Range("H2").Value = 6

Set-TextValue $q3Sheet.Range("B3") "004316"
Set-TextValue $q3Sheet.Range("C3") "前海开源沪港深裕鑫灵活配置混合A"
Set-TextValue $q3Sheet.Range("D3") "1.77"
Set-TextValue $q3Sheet.Range("E3") "70.17"
Set-TextValue $q3Sheet.Range("F3") "2.90"
Set-TextValue $q3Sheet.Range("G3") "0.0513"
$q3Sheet.Range("H3").Value = 6

Set-TextValue $q3Sheet.Range("B4") "001942"
Set-TextValue $q3Sheet.Range("C4") "前海开源沪港深汇鑫灵活配置混合A"
Set-TextValue $q3Sheet.Range("D4") "0.17"
Set-TextValue $q3Sheet.Range("E4") "87.24"
Set-TextValue $q3Sheet.Range("F4") "4.90"
Set-TextValue $q3Sheet.Range("G4") "0.0083"
$q3Sheet.Range("H4").Value = 3

Set-TextValue $q3Sheet.Range("B5") "001943"
Set-TextValue $q3Sheet.Range("C5") "前海开源沪港深汇鑫灵活配置混合C"
Set-TextValue $q3Sheet.Range("D5") "0.09"
Set-TextValue $q3Sheet.Range("E5") "87.24"
Set-TextValue $q3Sheet.Range("F5") "4.90"
Set-TextValue $q3Sheet.Range("G5") "0.0044"
$q3Sheet.Range("H5").Value = 3

# Reset the "@" text-number-format tweak back to the default/general style on
# every text cell we touched, using a format-only paste from a pristine,
# never-written cell far outside the used range.
$blank = $q3Sheet.Range("Z1")
$blank.Copy()
$q3Sheet.Range("B2:G5").PasteSpecial(-4122)

# =====================================================================
# Step 3: update the "总计" sheet - new 2022-Q3 row on top, everything else
# shifts down a row and gets re-numbered (column A is just a 0-based index).
# =====================================================================
$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = "2021-Q1"
$totalSheet.Range("C7").Value = 4
$totalSheet.Range("D7").Value = 0.18

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q3"
$totalSheet.Range("C6").Value = 2
$totalSheet.Range("D6").Value = 0.04

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q4"
$totalSheet.Range("C5").Value = 4
$totalSheet.Range("D5").Value = 0.08

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 3
$totalSheet.Range("D4").Value = 0.04

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.03

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.12

# The brand-new A7 cell needs the same "index column" style (s=2) the other
# A-column cells on this sheet already carry.
$totalSheet.Range("A6").Copy($totalSheet.Range("A7"))
$totalSheet.Range("A7").Value = 5
